$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.008.85'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.829.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '706.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.74'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.828.98'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.39'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.50'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.478.73'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.845.47'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.085.71'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.36'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '494.23'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.62'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.733'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.62'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.60'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.985.43'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.08%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.40'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.24'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.31'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.175'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.800.35'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.70%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.65%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.61%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.32'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.55%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '163.95'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("B47").Value = 'FLOKI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000310'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '429.68'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '48.95'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.76'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.296'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.56%  '
